$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71, shifting existing rows 71-120 down to 72-121
$ws.Rows("71").Insert()

# Populate the new row 71 with the data for this record
$ws.Range("A71").Value = 11
$ws.Range("B71").Value = "Vega Monumental Concepción"
$ws.Range("C71").Value = "Bíobío"
$ws.Range("D71").Value = 45062
$ws.Range("D71").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E71").Value = 8
$ws.Range("F71").Value = 100112012
$ws.Range("G71").Value = "Espinaca"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 50
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = 10400
$ws.Range("N71").Value = "$/cuna 10 kilos"
$ws.Range("O71").Value = "Región Metropolitana"
$ws.Range("P71").Value = 1040
$ws.Range("Q71").Value = 10
$ws.Range("R71").Value = "Hortaliza"
